$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new column header "MicroplasticImages" in Q1 (new shared string, new column)
$ws.Range("Q1").Value = "MicroplasticImages"

# Update the active selection to Q2 (as reflected in the sheet view)
$ws.Range("Q2").Select()
